# Add newly recorded capture-frequency counts to a few data rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-RowValues {
    param($ws, $rangeAddress, $vals)
    $arr = New-Object "object[,]" 1,$vals.Length
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $ws.Range($rangeAddress).Value = $arr
}

$rowValues = @(1, 0, 0, 0, 1, 0, 0, 0, 0, 0, 1)
Set-RowValues $ws "D51:N51" $rowValues

$rowValues = @(0, 0, 1, 0, 0, 5, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 1, 0, 1, 1, 1, 1, 2, 0, 0, 0, 1, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 1, 1, 1, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D55:CB55" $rowValues

$rowValues = @(0, 0, 1, 0, 0, 1, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 2, 0, 1, 0, 0, $null, 1, 0, 0, 0, 0, 0, 1, 0, 0, 0, 1, 1, 2, 0, 1, 0, 0, 0, 2, 3, 1, 2, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D60:CB60" $rowValues

$rowValues = @(1, 1, 0, 1, 0, 2, 2, 0, 2, 0, 0, 0, 1, 1, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws "D63:CB63" $rowValues

# Move the active selection to L49, as in the saved workbook state
$ws.Range("L49").Select()
